# Weekly update: a new week's observation is inserted at row 204, pushing the
# existing rows 204:237 down to 205:238 (dimension grows from R237 to R238).
# The new row 204 keeps the same market/category metadata as the old row 204
# but carries a fresh date (D) and volume (J) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("204:204").Insert()

$ws.Range("A204").Value = 10
$ws.Range("B204").Value = "Vega Modelo de Temuco"
$ws.Range("C204").Value = "La Araucanía"
$ws.Range("D204").Value = 44474
$ws.Range("E204").Value = 9
$ws.Range("F204").Value = 100112008
$ws.Range("G204").Value = "Coliflor"
$ws.Range("H204").Value = "Sin especificar"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 1800
$ws.Range("K204").Value = 800
$ws.Range("L204").Value = 800
$ws.Range("M204").Value = 800
$ws.Range("N204").Value = "$/unidad"
$ws.Range("O204").Value = "Región Metropolitana"
$ws.Range("P204").Value = 800
$ws.Range("Q204").Value = 1
$ws.Range("R204").Value = "Hortaliza"
